# Applies cryptos price/volume refresh per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain-text updates (percentages, and price strings that already contain
# multiple "." separators so Excel cannot parse them as a number).
$textUpdates = @(
    @('D2', '66.284.55'),
    @('E2', '  -0.82%  '),
    @('D3', '3.321.93'),
    @('E3', '  -0.71%  '),
    @('E4', '  -0.03%  '),
    @('E5', '  +2.74%  '),
    @('E6', '  +0.97%  '),
    @('E7', '  +8.25%  '),
    @('E8', '  +0.01%  '),
    @('E9', '  -2.71%  '),
    @('E10', '  +2.06%  '),
    @('E11', '  -0.31%  '),
    @('D12', '3.901.06'),
    @('E12', '  -0.69%  '),
    @('E13', '  -3.75%  '),
    @('D14', '66.312.06'),
    @('E14', '  -0.84%  '),
    @('E15', '  -3.36%  '),
    @('E16', '  -2.64%  '),
    @('D17', '3.277.56'),
    @('E17', '  -1.83%  '),
    @('E18', '  -2.08%  '),
    @('E19', '  -2.73%  '),
    @('E20', '  -3.24%  '),
    @('E21', '  -2.92%  '),
    @('E23', '  +0.14%  '),
    @('E24', '  +0.74%  '),
    @('D25', '3.461.76'),
    @('E26', '  -0.52%  '),
    @('E27', '  +7.93%  '),
    @('E28', '  -3.53%  '),
    @('E29', '  -1.43%  '),
    @('E30', '  -0.38%  '),
    @('E31', '  -0.51%  '),
    @('E32', '  -2.34%  '),
    @('E33', '  +0.09%  '),
    @('E34', '  -2.39%  '),
    @('E35', '  -3.12%  '),
    @('E36', '  -4.02%  '),
    @('E37', '  -1.02%  '),
    @('E38', '  -3.84%  '),
    @('D39', '2.880.84'),
    @('E39', '  +1.36%  '),
    @('E40', '  -2.41%  '),
    @('E41', '  -5.04%  '),
    @('E43', '  -2.56%  '),
    @('E44', '  -0.27%  '),
    @('E45', '  -1.24%  '),
    @('E46', '  -5.41%  '),
    @('E47', '  -2.61%  '),
    @('E48', '  -5.52%  '),
    @('E49', '  -2.89%  '),
    @('E50', '  -0.29%  '),
    @('E51', '  +5.02%  '),
)
foreach ($u in $textUpdates) {
    $ws.Range($u[0]).Value = $u[1]
}

# Price strings that look like plain numbers ("588.97", "0.998", ...) -- force
# them to stay literal text (matching the inline-string cells in the source
# workbook) instead of being auto-coerced to a floating point number, then
# restore the default "Normal" style so no stray number-format is left behind.
$numericLookingUpdates = @(
    @('D5', '588.97'),
    @('D6', '183.49'),
    @('D15', '26.25'),
    @('D18', '427.02'),
    @('D19', '5.53'),
    @('D21', '7.41'),
    @('D22', '71.94'),
    @('D24', '5.71'),
    @('D26', '0.517'),
    @('D29', '8.94'),
    @('D30', '0.998'),
    @('D32', '22.37'),
    @('D34', '5.18'),
    @('D37', '159.59'),
    @('D41', '26.47'),
    @('D43', '4.32'),
    @('D44', '40.13'),
    @('D45', '0.0665'),
    @('D46', '5.89'),
    @('D47', '2.31'),
    @('D49', '314.63'),
)
foreach ($u in $numericLookingUpdates) {
    $cell = $ws.Range($u[0])
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.Style = "Normal"
}

Write-Host "Applied $($textUpdates.Count + $numericLookingUpdates.Count) cell updates"
